# Conditional Aggregations Quiz Workbook - add COUNTIF/COUNTIFS/SUMIFS notes
# to column I, clear now-unused trailing rows/formatting, and tidy up the
# data columns' widths/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "How many employees" COUNTIF / COUNTIFS answers -----------
$ws.Range("I2").Formula  = '=COUNTIF($C:$C,"Cert4")'
$ws.Range("I3").Formula  = '=COUNTIF($C:$C,"Advanced Diploma")'
$ws.Range("I4").Formula  = '=COUNTIF(D:D,"Data")'
$ws.Range("I5").Formula  = '=COUNTIFS(E:E,"Data Scientist",C:C,"Advanced Diploma")'
$ws.Range("I6").Formula  = '=COUNTIF(B2:B20,"")'

# --- Add the "SUM of Salaries" SUMIFS answers ---------------------------
$ws.Range("I10").Formula = '=SUMIFS($F:$F,D:D,"Data")'
$ws.Range("I11").Formula = '=SUMIFS(F:F,C:C,"Masters Degree")'
$ws.Range("I12").Formula = '=SUMIFS(F:F,E:E,"Research Associate")'
$ws.Range("I13").Formula = '=SUMIFS(F:F,D:D,"Data",C:C,"Advanced Diploma")'
$ws.Range("I14").Formula = '=SUMIFS(F:F,E:E,"Data Scientist",C:C,"Bachelor Degree")'

# --- Remove now-empty leftover surname cells (no value, only formatting) -
$ws.Range("B8").Clear()
$ws.Range("B13").Clear()
$ws.Range("B15").Clear()
$ws.Range("B18").Clear()

# --- Strip the leftover cell formatting from the data table -------------
$ws.Range("A2:E20").Style = "Normal"

# --- Delete the unused, blank formatted rows below the data table -------
$ws.Range("A21:A51").EntireRow.Delete()

# --- Re-fit the data columns to their (now shorter) contents ------------
$ws.Columns.Item(1).ColumnWidth = 9.59
$ws.Columns.Item(2).ColumnWidth = 9.26
$ws.Columns.Item(3).ColumnWidth = 16.92
$ws.Columns.Item(4).ColumnWidth = 8.09
$ws.Columns.Item(5).ColumnWidth = 25.26
$ws.Columns.Item(6).ColumnWidth = 9.09

# --- Restore the selection to match the author's final cursor position --
$ws.Range("I14").Select()
